$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header renames (Chinese labels -> API field names) ---
$ws.Range("A1").Value = "isin"
$ws.Range("B1").Value = "lastPrice"
$ws.Range("C1").Value = "timestampLastPrice"
$ws.Range("D1").Value = "changeToPrevDayAbsolute"
$ws.Range("E1").Value = "changeToPrevDayInPercent"
$ws.Range("F1").Value = "closingPricePrevTradingDay"
$ws.Range("G1").Value = "mic"
$ws.Range("H1").Value = "lastPriceIndicator"
$ws.Range("I1").Value = "dayHigh"
$ws.Range("J1").Value = "dayLow"
$ws.Range("K1").Value = "priceFixings"
$ws.Range("L1").Value = "tradedInPercent"
$ws.Range("M1").Value = "tradingTimeEnd"
$ws.Range("N1").Value = "tradingTimeStart"
$ws.Range("O1").Value = "turnoverInEur"
$ws.Range("P1").Value = "turnoverInPieces"
$ws.Range("Q1").Value = "turnoverNominal"
$ws.Range("R1").Value = "weeks52High"
$ws.Range("S1").Value = "weeks52Low"
$ws.Range("T1").Value = "currency"

# New column U1 - copy formatting (bold/border/centered) from T1, then set its text
$ws.Range("T1").Copy($ws.Range("U1"))
$ws.Range("U1").Value = "minimumTradableUnit"

# --- Row 2: data values shift / update ---
$ws.Range("B2").Value = 90.44
$ws.Range("C2").Value = "2025-02-20T10:18:03+01:00"
$ws.Range("D2").Value = 0.37
$ws.Range("E2").Value = 0.41
$ws.Range("F2").Value = 90.44
# G2 (XFRA) unchanged

# H2 becomes a present-but-empty text cell (not a fully blank cell): force
# text entry via a lone leading apostrophe, then strip the resulting
# quote-prefix formatting by pasting A2's (default) format over it.
$ws.Range("H2").Formula = "'"
$ws.Range("A2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

$ws.Range("I2").Value = 90.44
$ws.Range("J2").Value = 90.44
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = $true
$ws.Range("M2").Value = "17:30:00"
$ws.Range("N2").Value = "08:00:00"
$ws.Range("O2").Value = 0

# P2 likewise becomes a present-but-empty text cell.
$ws.Range("P2").Formula = "'"
$ws.Range("A2").Copy()
$ws.Range("P2").PasteSpecial(-4122)

$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 101.29
$ws.Range("S2").Value = 87.34
$ws.Range("T2").Value = "{'originalValue': 'USD', 'translations': {'de': 'US-Dollar', 'en': 'U.S. dollar'}}"
$ws.Range("U2").Value = 200000
